$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking values remain stored as text, matching the
# "numberStoredAsText" convention already used throughout this sheet.
$ws.Range("C2:F12").NumberFormat = "@"

$ws.Range("C2").Value = "11"
$ws.Range("D2").Value = "10"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "0"

$ws.Range("C3").Value = "26"
$ws.Range("D3").Value = "15"
$ws.Range("E3").Value = "5"
$ws.Range("F3").Value = "0"

$ws.Range("C4").Value = "26"
$ws.Range("D4").Value = "19"
$ws.Range("E4").Value = "3"
$ws.Range("F4").Value = "0"

$ws.Range("C5").Value = "56"
$ws.Range("D5").Value = "39"
$ws.Range("E5").Value = "6"
$ws.Range("F5").Value = "1"

$ws.Range("C6").Value = "26"
$ws.Range("D6").Value = "20"
$ws.Range("E6").Value = "4"
$ws.Range("F6").Value = "0"

$ws.Range("C7").Value = "5"
$ws.Range("D7").Value = "9"
$ws.Range("E7").Value = "0"
$ws.Range("F7").Value = "0"

$ws.Range("C8").Value = "45"
$ws.Range("D8").Value = "25"
$ws.Range("E8").Value = "4"
$ws.Range("F8").Value = "3"

$ws.Range("C9").Value = "106"
$ws.Range("D9").Value = "50"
$ws.Range("E9").Value = "10"
$ws.Range("F9").Value = "7"

$ws.Range("C10").Value = "25"
$ws.Range("D10").Value = "18"
$ws.Range("E10").Value = "3"
$ws.Range("F10").Value = "0"

$ws.Range("C11").Value = "89"
$ws.Range("D11").Value = "60"
$ws.Range("E11").Value = "7"
$ws.Range("F11").Value = "4"

$ws.Range("C12").Value = "9"
$ws.Range("D12").Value = "6"
$ws.Range("E12").Value = "1"
$ws.Range("F12").Value = "0"

